$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '22.470.97'
Set-TextCell 'E2' '  +0.45%  '
Set-TextCell 'D3' '1.572.57'
Set-TextCell 'E3' '  +0.36%  '
Set-TextCell 'E4' '  -0.07%  '
Set-TextCell 'E5' '  +0.38%  '
Set-TextCell 'D6' '291.35'
Set-TextCell 'E6' '  -0.06%  '
Set-TextCell 'D7' '0.3704'
Set-TextCell 'E7' '  -1.61%  '
Set-TextCell 'D8' '49.93'
Set-TextCell 'E8' '  +1.40%  '
Set-TextCell 'D9' '0.3385'
Set-TextCell 'E9' '  -0.69%  '
Set-TextCell 'B10' 'Polygon'
Set-TextCell 'C10' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D10' '1.144'
Set-TextCell 'E10' '  +0.25%  '
Set-TextCell 'B11' 'Dogecoin'
Set-TextCell 'C11' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 'D11' '0.07550'
Set-TextCell 'E11' '  -0.79%  '
Set-TextCell 'D12' '1.001'
Set-TextCell 'E12' '  -0.15%  '
Set-TextCell 'D13' '21.30'
Set-TextCell 'E13' '  +1.21%  '
Set-TextCell 'D14' '6.040'
Set-TextCell 'E14' '  +1.04%  '
Set-TextCell 'D15' '6.968'
Set-TextCell 'E15' '  +0.31%  '
Set-TextCell 'D16' '1.570.85'
Set-TextCell 'E16' '  -0.86%  '
Set-TextCell 'E17' '  -1.14%  '
Set-TextCell 'D18' '90.71'
Set-TextCell 'E18' '  +0.87%  '
Set-TextCell 'D19' '0.06760'
Set-TextCell 'E19' '  +0.37%  '
Set-TextCell 'E20' '  +0.00%  '
Set-TextCell 'D21' '6.288'
Set-TextCell 'E21' '  +1.38%  '
Set-TextCell 'D22' '16.41'
Set-TextCell 'E22' '  -1.03%  '
Set-TextCell 'D23' '12.19'
Set-TextCell 'E23' '  +2.12%  '
Set-TextCell 'D24' '22.474.49'
Set-TextCell 'E24' '  +0.54%  '
Set-TextCell 'D25' '2.356'
Set-TextCell 'E25' '  -1.96%  '
Set-TextCell 'D26' '2.628'
Set-TextCell 'E26' '  -2.95%  '
Set-TextCell 'E27' '  -0.44%  '
Set-TextCell 'D28' '149.70'
Set-TextCell 'E28' '  +1.46%  '
Set-TextCell 'D29' '5.059'
Set-TextCell 'E29' '  +0.70%  '
Set-TextCell 'D30' '125.09'
Set-TextCell 'E30' '  -0.93%  '
Set-TextCell 'D31' '1.746.88'
Set-TextCell 'E31' '  +1.43%  '
Set-TextCell 'D32' '1.080'
Set-TextCell 'E32' '  +9.30%  '
Set-TextCell 'D33' '6.251'
Set-TextCell 'E33' '  +2.61%  '
Set-TextCell 'D34' '2.016'
Set-TextCell 'E34' '  +0.03%  '
Set-TextCell 'D35' '9.827'
Set-TextCell 'E35' '  -3.05%  '
Set-TextCell 'D36' '0.08365'
Set-TextCell 'E36' '  -1.32%  '
Set-TextCell 'E37' '  -1.30%  '
Set-TextCell 'D38' '0.2303'
Set-TextCell 'D39' '1.339'
Set-TextCell 'E39' '  -4.97%  '
Set-TextCell 'D40' '0.06555'
Set-TextCell 'D41' '5.460'
Set-TextCell 'E41' '  +1.05%  '
Set-TextCell 'E42' '  +0.37%  '
Set-TextCell 'D43' '0.6230'
Set-TextCell 'E43' '  -1.42%  '
Set-TextCell 'E44' '  -0.08%  '
Set-TextCell 'D45' '14.04'
Set-TextCell 'E45' '  +0.93%  '
Set-TextCell 'D46' '3.806'
Set-TextCell 'E46' '  -0.14%  '
Set-TextCell 'D47' '0.5866'
Set-TextCell 'E47' '  -1.05%  '
Set-TextCell 'D48' '129.26'
Set-TextCell 'E48' '  +3.77%  '
Set-TextCell 'D49' '2.076'
Set-TextCell 'E49' '  -0.26%  '
Set-TextCell 'D50' '1.221'
Set-TextCell 'E50' '  -4.20%  '
Set-TextCell 'D51' '0.07342'
Set-TextCell 'E51' '  +0.42%  '
